$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell value updates derived from the authoritative diff.
# Column D values that look like plain numbers get a leading apostrophe so
# Excel keeps storing them as text (matching the original inline-string type)
# instead of auto-converting them to numeric cells.

$ws.Range("D2").Value = "27.632.76"
$ws.Range("E2").Value = "  +0.00%  "
$ws.Range("D3").Value = "1.634.53"
$ws.Range("E3").Value = "  -0.16%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'212.26"
$ws.Range("E5").Value = "  -0.15%  "
$ws.Range("D6").Value = "'0.521"
$ws.Range("E6").Value = "  -0.51%  "
$ws.Range("E8").Value = "  +1.55%  "
$ws.Range("E9").Value = "  +2.76%  "
$ws.Range("D11").Value = "'0.0868"
$ws.Range("E11").Value = "  -2.62%  "
$ws.Range("D12").Value = "1.865.71"
$ws.Range("E12").Value = "  -0.17%  "
$ws.Range("D13").Value = "1.640.65"
$ws.Range("E13").Value = "  +0.35%  "
$ws.Range("E14").Value = "  +0.41%  "
$ws.Range("E15").Value = "  -0.81%  "
$ws.Range("D16").Value = "'65.26"
$ws.Range("E16").Value = "  +1.16%  "
$ws.Range("D17").Value = "27.618.18"
$ws.Range("E17").Value = "  -0.07%  "
$ws.Range("D18").Value = "'231.05"
$ws.Range("E18").Value = "  +0.87%  "
$ws.Range("E19").Value = "  -0.20%  "
$ws.Range("D20").Value = "'7.59"
$ws.Range("E20").Value = "  -2.01%  "
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("E22").Value = "  +5.66%  "
$ws.Range("D23").Value = "'4.36"
$ws.Range("E23").Value = "  +1.46%  "
$ws.Range("D24").Value = "'2.13"
$ws.Range("E24").Value = "  +6.82%  "
$ws.Range("D25").Value = "'149.56"
$ws.Range("E25").Value = "  -0.42%  "
$ws.Range("E26").Value = "  -0.92%  "
$ws.Range("E27").Value = "  +0.06%  "
$ws.Range("E28").Value = "  -0.14%  "
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("E30").Value = "  -0.21%  "
$ws.Range("E31").Value = "  -0.32%  "
$ws.Range("E32").Value = "  -0.46%  "
$ws.Range("D33").Value = "1.480.49"
$ws.Range("E33").Value = "  +1.80%  "
$ws.Range("E34").Value = "  -0.35%  "
$ws.Range("E35").Value = "  -1.23%  "
$ws.Range("E36").Value = "  -1.29%  "
$ws.Range("D37").Value = "'0.942"
$ws.Range("E37").Value = "  +4.70%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").Value = "'0.561"
$ws.Range("E38").Value = "  -0.18%  "
$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D39").Value = "'0.880"
$ws.Range("E39").Value = "  +0.65%  "
$ws.Range("E40").Value = "  +0.41%  "
$ws.Range("E41").Value = "  +2.33%  "
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("D43").Value = "'68.05"
$ws.Range("E43").Value = "  -2.49%  "
$ws.Range("D44").Value = "'2.47"
$ws.Range("E44").Value = "  +0.41%  "
$ws.Range("D45").Value = "'2.21"
$ws.Range("E45").Value = "  -1.16%  "
$ws.Range("D46").Value = "'5.35"
$ws.Range("E46").Value = "  -4.54%  "
$ws.Range("E47").Value = "  -0.22%  "
$ws.Range("E48").Value = "  +1.60%  "
$ws.Range("D49").Value = "'87.69"
$ws.Range("E49").Value = "  +1.56%  "
$ws.Range("E50").Value = "  -2.08%  "
$ws.Range("D51").Value = "'0.0993"
$ws.Range("E51").Value = "  +0.81%  "
